{"js": "// Replace the date line and each \"NNN\u00f7N=\" division prompt with its\n// updated value. Every <w:t> text run in the document is replaced\n// exactly once, and each original string is unique in the document,\n// so a plain search-and-replace (matchCase, not whole-word since the\n// strings contain punctuation) is safe and order-independent.\nconst replacements = [\n  [\"2025-02-21 Friday\", \"2025-02-22 Saturday\"],\n  [\"319\u00f79=\", \"217\u00f77=\"],\n  [\"487\u00f79=\", \"236\u00f79=\"],\n  [\"516\u00f75=\", \"839\u00f74=\"],\n  [\"698\u00f72=\", \"838\u00f76=\"],\n  [\"200\u00f77=\", \"299\u00f75=\"],\n  [\"836\u00f79=\", \"632\u00f78=\"],\n  [\"360\u00f79=\", \"540\u00f73=\"],\n  [\"694\u00f78=\", \"103\u00f77=\"],\n  [\"711\u00f78=\", \"467\u00f77=\"],\n  [\"417\u00f75=\", \"170\u00f75=\"],\n  [\"390\u00f79=\", \"830\u00f79=\"],\n  [\"382\u00f72=\", \"845\u00f73=\"],\n  [\"331\u00f72=\", \"410\u00f79=\"],\n  [\"571\u00f75=\", \"575\u00f72=\"],\n  [\"595\u00f79=\", \"352\u00f72=\"],\n  [\"131\u00f72=\", \"753\u00f79=\"],\n  [\"433\u00f77=\", \"417\u00f72=\"],\n  [\"676\u00f78=\", \"511\u00f79=\"],\n  [\"566\u00f77=\", \"767\u00f79=\"],\n  [\"725\u00f79=\", \"318\u00f79=\"],\n  [\"203\u00f74=\", \"820\u00f76=\"],\n  [\"287\u00f76=\", \"803\u00f77=\"],\n  [\"616\u00f78=\", \"916\u00f79=\"],\n  [\"482\u00f72=\", \"322\u00f72=\"],\n  [\"657\u00f77=\", \"368\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each \"NNN\u00f7N=\" division prompt with its\n# updated value. Every text run in the document is replaced exactly\n# once, and each original string is unique in the document, so a plain\n# Find/Replace-all per pair is safe and order-independent.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-02-21 Friday\", \"2025-02-22 Saturday\"),\n    @(\"319\u00f79=\", \"217\u00f77=\"),\n    @(\"487\u00f79=\", \"236\u00f79=\"),\n    @(\"516\u00f75=\", \"839\u00f74=\"),\n    @(\"698\u00f72=\", \"838\u00f76=\"),\n    @(\"200\u00f77=\", \"299\u00f75=\"),\n    @(\"836\u00f79=\", \"632\u00f78=\"),\n    @(\"360\u00f79=\", \"540\u00f73=\"),\n    @(\"694\u00f78=\", \"103\u00f77=\"),\n    @(\"711\u00f78=\", \"467\u00f77=\"),\n    @(\"417\u00f75=\", \"170\u00f75=\"),\n    @(\"390\u00f79=\", \"830\u00f79=\"),\n    @(\"382\u00f72=\", \"845\u00f73=\"),\n    @(\"331\u00f72=\", \"410\u00f79=\"),\n    @(\"571\u00f75=\", \"575\u00f72=\"),\n    @(\"595\u00f79=\", \"352\u00f72=\"),\n    @(\"131\u00f72=\", \"753\u00f79=\"),\n    @(\"433\u00f77=\", \"417\u00f72=\"),\n    @(\"676\u00f78=\", \"511\u00f79=\"),\n    @(\"566\u00f77=\", \"767\u00f79=\"),\n    @(\"725\u00f79=\", \"318\u00f79=\"),\n    @(\"203\u00f74=\", \"820\u00f76=\"),\n    @(\"287\u00f76=\", \"803\u00f77=\"),\n    @(\"616\u00f78=\", \"916\u00f79=\"),\n    @(\"482\u00f72=\", \"322\u00f72=\"),\n    @(\"657\u00f77=\", \"368\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
